$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.242386666666667
$ws.Range("H2").Value = 6.72716
$ws.Range("I2").Value = 0.04442500453715972
$ws.Range("J2").Value = 0.04442500453715972
$ws.Range("M2").Value = 10.92359866666667
$ws.Range("N2").Value = 32.770796
$ws.Range("O2").Value = 0.2236009040380497
$ws.Range("P2").Value = 0.2236009040380497
$ws.Range("Q2").Value = 24.49493200215111
$ws.Range("R2").Value = 220.45438801936
$ws.Range("S2").Value = 0.009933471176403374
$ws.Range("T2").Value = 0.00993347117640337

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.242386666666667
$ws.Range("H3").Value = 6.72716
$ws.Range("I3").Value = 0.04442500453715972
$ws.Range("J3").Value = 0.04442500453715972
$ws.Range("O3").Value = 0.4261214970992155
$ws.Range("P3").Value = 0.4261214970992155
$ws.Range("Q3").Value = 46.68056750935111
$ws.Range("R3").Value = 420.1251075841599
$ws.Range("S3").Value = 0.01893044944201394
$ws.Range("T3").Value = 0.01893044944201394

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.242386666666667
$ws.Range("H4").Value = 6.72716
$ws.Range("I4").Value = 0.04442500453715972
$ws.Range("J4").Value = 0.04442500453715972
$ws.Range("M4").Value = 13.06524766666667
$ws.Range("N4").Value = 39.195743
$ws.Range("O4").Value = 0.2674394472823625
$ws.Range("P4").Value = 0.2674394472823625
$ws.Range("Q4").Value = 29.29733716443111
$ws.Range("R4").Value = 263.67603447988
$ws.Range("S4").Value = 0.01188099865893444
$ws.Range("T4").Value = 0.01188099865893444

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.242386666666667
$ws.Range("H5").Value = 6.72716
$ws.Range("I5").Value = 0.04442500453715972
$ws.Range("J5").Value = 0.04442500453715972
$ws.Range("M5").Value = 4.046901
$ws.Range("N5").Value = 12.140703
$ws.Range("O5").Value = 0.0828381515803724
$ws.Range("P5").Value = 0.0828381515803724
$ws.Range("Q5").Value = 9.074716843719999
$ws.Range("R5").Value = 81.67245159347999
$ws.Range("S5").Value = 0.003680085259807969
$ws.Range("T5").Value = 0.003680085259807968

# Row 6
$ws.Range("I6").Value = 0.4052409520727612
$ws.Range("J6").Value = 0.4052409520727612
$ws.Range("M6").Value = 10.92359866666667
$ws.Range("N6").Value = 32.770796
$ws.Range("O6").Value = 0.2236009040380497
$ws.Range("P6").Value = 0.2236009040380497
$ws.Range("Q6").Value = 223.4405976752636
$ws.Range("R6").Value = 2010.965379077372
$ws.Range("S6").Value = 0.09061224323670938
$ws.Range("T6").Value = 0.09061224323670936

# Row 7
$ws.Range("I7").Value = 0.4052409520727612
$ws.Range("J7").Value = 0.4052409520727612
$ws.Range("O7").Value = 0.4261214970992155
$ws.Range("P7").Value = 0.4261214970992155
$ws.Range("S7").Value = 0.1726818811831565
$ws.Range("T7").Value = 0.1726818811831564

# Row 8
$ws.Range("I8").Value = 0.4052409520727612
$ws.Range("J8").Value = 0.4052409520727612
$ws.Range("M8").Value = 13.06524766666667
$ws.Range("N8").Value = 39.195743
$ws.Range("O8").Value = 0.2674394472823625
$ws.Range("P8").Value = 0.2674394472823625
$ws.Range("Q8").Value = 267.2477117200945
$ws.Range("R8").Value = 2405.229405480851
$ws.Range("S8").Value = 0.1083774162385176
$ws.Range("T8").Value = 0.1083774162385176

# Row 9
$ws.Range("I9").Value = 0.4052409520727612
$ws.Range("J9").Value = 0.4052409520727612
$ws.Range("M9").Value = 4.046901
$ws.Range("N9").Value = 12.140703
$ws.Range("O9").Value = 0.0828381515803724
$ws.Range("P9").Value = 0.0828381515803724
$ws.Range("Q9").Value = 82.77876236261899
$ws.Range("R9").Value = 745.008861263571
$ws.Range("S9").Value = 0.03356941141437782
$ws.Range("T9").Value = 0.03356941141437782

# Row 10
$ws.Range("G10").Value = 27.778539
$ws.Range("H10").Value = 83.335617
$ws.Range("I10").Value = 0.5503340433900792
$ws.Range("J10").Value = 0.5503340433900791
$ws.Range("M10").Value = 10.92359866666667
$ws.Range("N10").Value = 32.770796
$ws.Range("O10").Value = 0.2236009040380497
$ws.Range("P10").Value = 0.2236009040380497
$ws.Range("Q10").Value = 303.4416115823481
$ws.Range("R10").Value = 2730.974504241132
$ws.Range("S10").Value = 0.123055189624937
$ws.Range("T10").Value = 0.1230551896249369

# Row 11
$ws.Range("G11").Value = 27.778539
$ws.Range("H11").Value = 83.335617
$ws.Range("I11").Value = 0.5503340433900792
$ws.Range("J11").Value = 0.5503340433900791
$ws.Range("O11").Value = 0.4261214970992155
$ws.Range("P11").Value = 0.4261214970992155
$ws.Range("Q11").Value = 578.275809598988
$ws.Range("R11").Value = 5204.482286390892
$ws.Range("S11").Value = 0.2345091664740452
$ws.Range("T11").Value = 0.2345091664740451

# Row 12
$ws.Range("G12").Value = 27.778539
$ws.Range("H12").Value = 83.335617
$ws.Range("I12").Value = 0.5503340433900792
$ws.Range("J12").Value = 0.5503340433900791
$ws.Range("M12").Value = 13.06524766666667
$ws.Range("N12").Value = 39.195743
$ws.Range("O12").Value = 0.2674394472823625
$ws.Range("P12").Value = 0.2674394472823625
$ws.Range("Q12").Value = 362.933491853159
$ws.Range("R12").Value = 3266.401426678431
$ws.Range("S12").Value = 0.1471810323849105
$ws.Range("T12").Value = 0.1471810323849105

# Row 13
$ws.Range("G13").Value = 27.778539
$ws.Range("H13").Value = 83.335617
$ws.Range("I13").Value = 0.5503340433900792
$ws.Range("J13").Value = 0.5503340433900791
$ws.Range("M13").Value = 4.046901
$ws.Range("N13").Value = 12.140703
$ws.Range("O13").Value = 0.0828381515803724
$ws.Range("P13").Value = 0.0828381515803724
$ws.Range("Q13").Value = 112.416997257639
$ws.Range("R13").Value = 1011.752975318751
$ws.Range("S13").Value = 0.04558865490618663
$ws.Range("T13").Value = 0.04558865490618661
